# Regenerate the sval data (filtered save games) for sauer_matt 2024 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("0.0008583669626518464", "0.002777888934908601", "0.8054896365839992", "0.496779210170732", "1.305905102652292")
    3 = @("0.04763786555579896",   "0.04240448674262143",  "0.8054896365839992", "0.496779210170732", "1.392311199053152")
    4 = @("0.6753301551942219",    "0.3127903958511391",   "0.8054896365839992", "0.496779210170732", "2.290389397800092")
    5 = @("0.0008583669626518464", "6.633126561350622e-07","0.1575252929769615", "0.496779210170732", "0.6551635334230015")
    6 = @("1.459612070389937",     "1.667794583268128",    "0.8054896365839992", "0.496779210170732", "4.429675500412797")
    7 = @("1.459612070389937",     "1.667794583268128",    "0.1575252929769615", "0.496779210170732", "3.781711156805759")
    8 = @("0.6753301551942219",    "1.667794583268128",    "26.21740644021617",  "0.496779210170732", "29.05731038884925")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = [double]$vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = [double]$vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = [double]$vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = [double]$vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = [double]$vals[4]  # G - sum
}
